# Generate Report for Handoff
# Updates the "dab8b466-2b5c-4135-9ca7-428a55f8e560" row's latest
# handoff/handback timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-27-12 00:27:45"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-12 00:27:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-12 00:27:45"
